# Update the three Persian description cells (SubjectDescriptionTranslated column)
# to match the revised wording from the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - "Dimension" subject: description translated (column E)
$ws.Range("E6").Value = "ابعاد باید با ابزار دارای گواهی کالیبراسیون  و بر اساس نقشه و استاندارد مربوطه اندازه گیری شود. تمام اندازه گیری های واقعی باید در جدول مناسب در گزارش بازرسی ذکر شود."

# Row 9 - "Hydrostatic test" subject: description translated (column E)
$ws.Range("E9").Value = "تست هیدرواستاتیک/فشار باید طبق دستورالعمل تایید شده و استاندارد مربوطه انجام شود. در گزارش باید شرایط محیطی از قبیل fitting و دمای آب در نظر گرفته شود. قبل از آزمایش، تمام گواهی های کالیبراسیون برای گیج ها و آنالیز آب باید بررسی شوند. (تاریخ انقضا و محتوای CL/PH مهم هستند)."

# Row 11 - "Document review" subject: description translated (column E)
$ws.Range("E11").Value = "بررسی سند: MTC باید بررسی شود و تایید شود که آنها با لوله های بازرسی شده مرتبط هستند. بازرس باید در گزارش ذکر کند که گواهینامه صهدر شده اصل بوده یا کپی و همچنین نوع آن: 2.1،2.2،3.1 یا 3.2 همچنین گواهی کیفیت مواد صادر شده توسط کارخانه/فروشنده در لیست فروشنده تایید شده است یا خیر. کلیه پارامترها مانند آنالیزهای شیمیایی، خواص مکانیکی و غیره باید با استاندارد مربوطه مقایسه شوند. بازرس فقط در صورتی که همه موارد رضایت بخش بود، بازرس تنها مجاز به استفاده از مهر review یاذکر این در امضا خود است."

# Update the view/selection state: scroll so column B is at the left edge
# and the active selection moves from I21 to F21.
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("F21").Select()
